$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2026-01-04 01:44:08"

# Update the "取得日時" (retrieved at) timestamp for the existing rows (2-6)
$ws.Range("A2").Value = $timestamp
$ws.Range("A3").Value = $timestamp
$ws.Range("A4").Value = $timestamp
$ws.Range("A5").Value = $timestamp
$ws.Range("A6").Value = $timestamp

# Append the new row (7) with the newly scraped listing
$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "【急募】簡単なHP作成とAWS構築をしてくれる方募集"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5457524"
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5457524")
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("G7").Value = 18
